$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44893
$ws.Cells.Item(2, 10).Value = 3300
$ws.Cells.Item(2, 11).Value = 1200
$ws.Cells.Item(2, 12).Value = 1300
$ws.Cells.Item(2, 13).Value = 1261
$ws.Cells.Item(2, 16).Value = 1261

$ws.Cells.Item(3, 4).Value = 45210
$ws.Cells.Item(3, 10).Value = 550
$ws.Cells.Item(3, 11).Value = 1500
$ws.Cells.Item(3, 12).Value = 1600
$ws.Cells.Item(3, 13).Value = 1536
$ws.Cells.Item(3, 16).Value = 1536

$ws.Cells.Item(4, 4).Value = 45203
$ws.Cells.Item(4, 10).Value = 800
$ws.Cells.Item(4, 11).Value = 1800
$ws.Cells.Item(4, 12).Value = 2000
$ws.Cells.Item(4, 13).Value = 1900
$ws.Cells.Item(4, 16).Value = 1900

$ws.Cells.Item(5, 4).Value = 45205
$ws.Cells.Item(5, 10).Value = 3500
$ws.Cells.Item(5, 11).Value = 1400
$ws.Cells.Item(5, 12).Value = 1500
$ws.Cells.Item(5, 13).Value = 1457
$ws.Cells.Item(5, 16).Value = 1457

$ws.Cells.Item(6, 4).Value = 44175
$ws.Cells.Item(6, 10).Value = 1400
$ws.Cells.Item(6, 11).Value = 1900
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 1950
$ws.Cells.Item(6, 16).Value = 1950

$ws.Cells.Item(7, 4).Value = 45204
$ws.Cells.Item(7, 10).Value = 1200
$ws.Cells.Item(7, 11).Value = 1600
$ws.Cells.Item(7, 12).Value = 1700
$ws.Cells.Item(7, 13).Value = 1650
$ws.Cells.Item(7, 16).Value = 1650

$ws.Cells.Item(8, 4).Value = 45062
$ws.Cells.Item(8, 10).Value = 1700
$ws.Cells.Item(8, 11).Value = 2800
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = 2900
$ws.Cells.Item(8, 16).Value = 2900

$ws.Cells.Item(9, 4).Value = 44895
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 1200
$ws.Cells.Item(9, 12).Value = 1300
$ws.Cells.Item(9, 13).Value = 1255
$ws.Cells.Item(9, 16).Value = 1255

$ws.Cells.Item(10, 4).Value = 44210
$ws.Cells.Item(10, 10).Value = 1450
$ws.Cells.Item(10, 11).Value = 1600
$ws.Cells.Item(10, 12).Value = 1700
$ws.Cells.Item(10, 13).Value = 1650
$ws.Cells.Item(10, 16).Value = 1650

$ws.Cells.Item(11, 4).Value = 45212
$ws.Cells.Item(11, 10).Value = 750
$ws.Cells.Item(11, 11).Value = 1400
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 13).Value = 1440
$ws.Cells.Item(11, 16).Value = 1440

$ws.Cells.Item(12, 4).Value = 44200
$ws.Cells.Item(12, 10).Value = 1500
$ws.Cells.Item(12, 11).Value = 1400
$ws.Cells.Item(12, 12).Value = 1500
$ws.Cells.Item(12, 13).Value = 1450
$ws.Cells.Item(12, 16).Value = 1450

$ws.Cells.Item(13, 4).Value = 44537
$ws.Cells.Item(13, 10).Value = 800
$ws.Cells.Item(13, 11).Value = 1300
$ws.Cells.Item(13, 12).Value = 1400
$ws.Cells.Item(13, 13).Value = 1350
$ws.Cells.Item(13, 16).Value = 1350

$ws.Cells.Item(14, 4).Value = 45132
$ws.Cells.Item(14, 10).Value = 170
$ws.Cells.Item(14, 11).Value = 2200
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 13).Value = 2359
$ws.Cells.Item(14, 16).Value = 2359

$ws.Cells.Item(15, 4).Value = 44907
$ws.Cells.Item(15, 10).Value = 2300
$ws.Cells.Item(15, 11).Value = 900
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 13).Value = 952
$ws.Cells.Item(15, 16).Value = 952

$ws.Cells.Item(16, 4).Value = 44883
$ws.Cells.Item(16, 10).Value = 290
$ws.Cells.Item(16, 11).Value = 1400
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = 1434
$ws.Cells.Item(16, 16).Value = 1434

$ws.Cells.Item(17, 4).Value = 44638
$ws.Cells.Item(17, 10).Value = 800
$ws.Cells.Item(17, 11).Value = 2500
$ws.Cells.Item(17, 12).Value = 2800
$ws.Cells.Item(17, 13).Value = 2650
$ws.Cells.Item(17, 16).Value = 2650
